$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.060.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +4.90%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.087.84'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.06%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.23'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.94'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +3.17%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.075.22'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.88%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.29%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.07%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +8.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.468'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.56%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000241'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.91%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.23'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.51%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.599.37'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.27'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.56%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.087.38'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '61.984.23'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '448.89'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.95'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.51'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +6.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.85'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.09'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.26'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.92%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.98%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.13'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.78'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +10.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.112'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +13.45%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.63'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.04'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.23%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0801'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.06'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.19'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +5.77%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.39'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.00'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +9.71%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.84'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '432.65'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +8.71%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.790.90'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.60%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.21%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +7.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '35.53'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +10.38%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.52'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.01%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.20%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.02'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.55%  '
